# Kammari_LabExam03Grading.xlsx - grading rubric updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 20 (S.No 12 - addProduct() method): grading comment + points changed
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = "(-2)For if the customer does not exists a new linked list should be initialized"

# Row 34 (S.No 18 - 100% passing of all the test cases): points + grading comment changed
$ws.Range("E34").Value = 4
$ws.Range("F34").Value = "(-3)I have changed your addProduct() code and run the test cases then 3 test cases failed but I didn’t deducted any points for remaining test cases.. "

# Update the active selection / scrolled view to match the reviewer's final position
$ws.Range("E31").Select()
